$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab "Scanner" -> "Morgue"
$ws.Name = "Morgue"

# Header row: "Type" -> "Number"
$ws.Range("E1").Value = "Number"

# Data row 2: date/time became invalid ("NaN/.." placeholders) and the
# "Type" column now holds a raw epoch-millis timestamp (kept as text).
$ws.Range("C2").Value = "NaN/NaN/NaN"
$ws.Range("D2").Value = "NaN:NaN:NaN"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1746364263817"
$ws.Range("E2").Style = "Normal"
